$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "... C#.Net, ADO.Net, HTML, CSS" -> "... C#.Net, ADO.Net, EntityFramework, HTML, CSS"
# The run holding " ADO.Net," needs to become three runs:
#   " ADO.Net"  |  ", EntityFramework"  |  ","
# ---------------------------------------------------------------------------

# Narrow down to the unique paragraph fragment first so the later, shorter
# search (" ADO.Net,") cannot match any of the several other "ADO.Net"
# occurrences elsewhere in the document.
$anchor = $d.Content
$anchorFound = $anchor.Find.Execute("C#.Net, ADO.Net, HTML, CSS", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)

if ($anchorFound) {
    $scoped = $d.Range($anchor.Start, $anchor.End)
    $hit = $scoped.Find.Execute(" ADO.Net,", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)

    if ($hit) {
        $newText = " ADO.Net, EntityFramework,"

        # Replace the run's text in place (keeps its existing formatting).
        $scoped.Text = $newText

        $s = $scoped.Start
        $e = $scoped.End

        $part1Len = " ADO.Net".Length
        $part2Len = ", EntityFramework".Length

        $r1 = $d.Range($s, $s + $part1Len)
        $r2 = $d.Range($s + $part1Len, $s + $part1Len + $part2Len)
        $r3 = $d.Range($s + $part1Len + $part2Len, $e)

        # Re-stamping (unchanged) character formatting on each fragment forces
        # the engine to keep them as separate runs instead of re-merging them
        # with neighbouring runs that happen to share identical rPr.
        foreach ($part in @($r1, $r2, $r3)) {
            $color = $part.Font.Color
            $part.Font.Color = 1
            $part.Font.Color = $color
        }
    }
}

# ---------------------------------------------------------------------------
# Edit 2: mark the run that carries the page-break + picture as NoProof.
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $shape = $d.InlineShapes.Item(1)
    $shapeRange = $shape.Range
    $shapeRange.NoProofing = $true
}
